$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()
$ws.Range("H135").Value = 41667450
$ws.Range("I135").Value = 16129833
$ws.Range("J135").Value = 200000660
$ws.Range("K135").Value = 145168497
$ws.Range("L135").Value = 1800005940
$ws.Range("M135").Value = -145165962
$ws.Range("N135").Value = -1800011010
$ws.Range("H138").Value = 4332.242
$ws.Range("I138").Value = 1378.9354
$ws.Range("J138").Value = 50108.5
$ws.Range("K138").Value = 4136.8062
$ws.Range("L138").Value = 150325.5
$ws.Range("M138").Value = 1003.1938
$ws.Range("N138").Value = -160605.5
$ws.Range("H141").Value = 3575.8333
$ws.Range("I141").Value = 2890.3845
$ws.Range("K141").Value = 8671.1535
$ws.Range("M141").Value = -3491.1535

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20641.3
$ws.Range("I32").Value = 21485.254
$ws.Range("J32").Value = 14839.125
$ws.Range("K32").Value = 21485.254
$ws.Range("L32").Value = 14839.125
$ws.Range("M32").Value = -21198.254
$ws.Range("N32").Value = -15413.125
$ws.Range("H122").Value = 12501449
$ws.Range("I122").Value = 1610.2222
$ws.Range("J122").Value = 125000000
$ws.Range("K122").Value = 4830.6666
$ws.Range("L122").Value = 375000000
$ws.Range("M122").Value = -2380.6666
$ws.Range("N122").Value = -375004900

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1970.122
$ws.Range("I134").Value = 2029.4231
$ws.Range("J134").Value = 1867.3334
$ws.Range("K134").Value = 6088.2693
$ws.Range("L134").Value = 5602.0002
$ws.Range("M134").Value = -3553.2693
$ws.Range("N134").Value = -10672.0002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 204.75
$ws.Range("I22").Value = 204.75
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 204.75
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 145.25
$ws.Range("N22").ClearContents()
$ws.Range("H58").Value = 1542591.4
$ws.Range("I58").Value = 2332187.5
$ws.Range("K58").Value = 2332187.5
$ws.Range("M58").Value = -2331984.5
$ws.Range("H64").Value = 40135.5
$ws.Range("J64").Value = 40135.5
$ws.Range("L64").Value = 40135.5
$ws.Range("N64").Value = -40631.5
$ws.Range("H67").Value = 40135.5
$ws.Range("J67").Value = 40135.5
$ws.Range("L67").Value = 40135.5
$ws.Range("N67").Value = -41851.5
$ws.Range("H122").Value = 9439.727999999999
$ws.Range("I122").Value = 6218.154
$ws.Range("J122").Value = 14093.111
$ws.Range("K122").Value = 18654.462
$ws.Range("L122").Value = 42279.333
$ws.Range("M122").Value = -16204.462
$ws.Range("N122").Value = -47179.333
$ws.Range("H123").Value = 53685
$ws.Range("J123").Value = 53685
$ws.Range("L123").Value = 53685
$ws.Range("N123").Value = -63485
$ws.Range("H132").Value = 1920.7358
$ws.Range("I132").Value = 1618.8529
$ws.Range("J132").Value = 2460.9473
$ws.Range("K132").Value = 4856.5587
$ws.Range("L132").Value = 7382.841899999999
$ws.Range("M132").Value = -2326.5587
$ws.Range("N132").Value = -12442.8419
$ws.Range("H136").Value = 1542591.4
$ws.Range("I136").Value = 2332187.5
$ws.Range("K136").Value = 6996562.5
$ws.Range("M136").Value = -6994012.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 717.1923
$ws.Range("J113").Value = 730.3077
$ws.Range("L113").Value = 2190.9231
$ws.Range("N113").Value = -6530.9231

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 19800
$ws.Range("J21").Value = 19800
$ws.Range("L21").Value = 19800
$ws.Range("N21").Value = -20146
$ws.Range("H24").Value = 1119888.2
$ws.Range("I24").Value = 20000000
$ws.Range("J24").Value = 9293.471
$ws.Range("K24").Value = 20000000
$ws.Range("L24").Value = 9293.471
$ws.Range("M24").Value = -19999827
$ws.Range("N24").Value = -9639.471
$ws.Range("H30").Value = 19800
$ws.Range("J30").Value = 19800
$ws.Range("L30").Value = 19800
$ws.Range("N30").Value = -20010
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()
$ws.Range("H64").Value = 38271
$ws.Range("J64").Value = 38271
$ws.Range("L64").Value = 38271
$ws.Range("N64").Value = -38767
$ws.Range("H67").Value = 38271
$ws.Range("J67").Value = 38271
$ws.Range("L67").Value = 38271
$ws.Range("N67").Value = -39987
$ws.Range("H107").Value = 407.14285
$ws.Range("I107").Value = 175.77777
$ws.Range("J107").Value = 823.6
$ws.Range("K107").Value = 175.77777
$ws.Range("L107").Value = 823.6
$ws.Range("M107").Value = 1744.22223
$ws.Range("N107").Value = -4663.6
$ws.Range("H122").Value = 9955.308000000001
$ws.Range("I122").Value = 15471.286
$ws.Range("J122").Value = 3520
$ws.Range("K122").Value = 46413.858
$ws.Range("L122").Value = 10560
$ws.Range("M122").Value = -43963.858
$ws.Range("N122").Value = -15460
$ws.Range("H132").Value = 2550.2
$ws.Range("I132").Value = 2359.84
$ws.Range("J132").Value = 3502
$ws.Range("K132").Value = 7079.52
$ws.Range("L132").Value = 10506
$ws.Range("M132").Value = -4549.52
$ws.Range("N132").Value = -15566

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 8333.333000000001
$ws.Range("I14").Value = 5000
$ws.Range("J14").Value = 10000
$ws.Range("K14").Value = 5000
$ws.Range("L14").Value = 10000
$ws.Range("M14").Value = -4828
$ws.Range("N14").Value = -10344
$ws.Range("H22").Value = 960.3333
$ws.Range("I22").Value = 890.7143
$ws.Range("J22").Value = 1004.63635
$ws.Range("K22").Value = 890.7143
$ws.Range("L22").Value = 1004.63635
$ws.Range("M22").Value = -595.7143
$ws.Range("N22").Value = -1594.63635
$ws.Range("H27").Value = 960.3333
$ws.Range("I27").Value = 890.7143
$ws.Range("J27").Value = 1004.63635
$ws.Range("K27").Value = 890.7143
$ws.Range("L27").Value = 1004.63635
$ws.Range("M27").Value = -783.7143
$ws.Range("N27").Value = -1218.63635
$ws.Range("H46").Value = 822.2143
$ws.Range("I46").Value = 801.5714
$ws.Range("J46").Value = 842.8570999999999
$ws.Range("K46").Value = 801.5714
$ws.Range("L46").Value = 842.8570999999999
$ws.Range("M46").Value = -613.5714
$ws.Range("N46").Value = -1218.8571
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H100").Value = 5322.154
$ws.Range("I100").Value = 2636
$ws.Range("K100").Value = 2636
$ws.Range("M100").Value = -2095
$ws.Range("H122").Value = 6802.327
$ws.Range("I122").Value = 6524.086
$ws.Range("J122").Value = 7289.25
$ws.Range("K122").Value = 19572.258
$ws.Range("L122").Value = 21867.75
$ws.Range("M122").Value = -17122.258
$ws.Range("N122").Value = -26767.75
$ws.Range("H132").Value = 10062.5
$ws.Range("I132").Value = 15820.923
$ws.Range("J132").Value = 4304.077
$ws.Range("K132").Value = 47462.769
$ws.Range("L132").Value = 12912.231
$ws.Range("M132").Value = -44932.769
$ws.Range("N132").Value = -17972.231
$ws.Range("H136").Value = 4172.396
$ws.Range("I136").Value = 2205.0605
$ws.Range("K136").Value = 6615.181500000001
$ws.Range("M136").Value = -4065.181500000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 19800
$ws.Range("J26").Value = 19800
$ws.Range("L26").Value = 19800
$ws.Range("N26").Value = -20386
$ws.Range("H37").Value = 29610
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 29610
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 29610
$ws.Range("M37").ClearContents()
$ws.Range("N37").Value = -30016
$ws.Range("H63").Value = 40249
$ws.Range("J63").Value = 40249
$ws.Range("L63").Value = 40249
$ws.Range("N63").Value = -41497
$ws.Range("H66").Value = 40249
$ws.Range("J66").Value = 40249
$ws.Range("L66").Value = 120747
$ws.Range("N66").Value = -126987
$ws.Range("H132").Value = 1555.3877
$ws.Range("I132").Value = 763.7105
$ws.Range("J132").Value = 4290.273
$ws.Range("K132").Value = 2291.1315
$ws.Range("L132").Value = 12870.819
$ws.Range("M132").Value = 238.8685
$ws.Range("N132").Value = -17930.819
$ws.Range("H136").Value = 6220.283
$ws.Range("I136").Value = 4373.029
$ws.Range("K136").Value = 13119.087
$ws.Range("M136").Value = -10569.087
